$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: phone number was stored as text "09876543" with a leading zero.
# Convert it to a plain number (9876543), keeping the existing points (120).
$ws.Cells.Item(36, 1).Value = 9876543

# Add new row 37: re-add the original text phone number "09876543" (with
# leading zero, as text) with its points reset to 0.
$ws.Cells.Item(37, 1).NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = "09876543"
$ws.Cells.Item(37, 1).Style = "Normal"
$ws.Cells.Item(37, 2).NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = ""
$ws.Cells.Item(37, 2).Style = "Normal"
$ws.Cells.Item(37, 3).Value = 0
